$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2378.06569508192
$ws.Range("I2").Value = 752.065695081923
$ws.Range("B3").Value = 2343.66497734949
$ws.Range("E3").Value = 2649.05235760662
$ws.Range("F3").Value = 2831.65681407097
$ws.Range("I3").Value = 1298.66497734949
$ws.Range("B4").Value = 3103.1995245717
$ws.Range("C4").Value = 2579.09699842962
$ws.Range("D4").Value = 2518.22373248302
$ws.Range("I4").Value = 1462.1995245717
$ws.Range("B5").Value = 2907.35067905523
$ws.Range("E5").Value = 3294.83628293049
$ws.Range("I5").Value = 1175.35067905523
$ws.Range("B6").Value = 2630.29632699167
$ws.Range("E6").Value = 3168.26125267739
$ws.Range("F6").Value = 3238.77355413596
$ws.Range("I6").Value = 1045.29632699167
$ws.Range("B7").Value = 2269.50489783241
$ws.Range("E7").Value = 2664.1207550185
$ws.Range("F7").Value = 2847.89966836922
$ws.Range("I7").Value = 561.504897832414
$ws.Range("B8").Value = 2300.54913248253
$ws.Range("E8").Value = 2740.24638034679
$ws.Range("I8").Value = 541.549132482534
$ws.Range("B9").Value = 2271.11056592564
$ws.Range("I9").Value = 512.110565925638
$ws.Range("B10").Value = 2022.94335389799
$ws.Range("E10").Value = 2399.62059869856
$ws.Range("I10").Value = 211.943353897987
$ws.Range("B11").Value = 1867.69979525656
$ws.Range("E11").Value = 2304.32527005159
$ws.Range("F11").Value = 2429.22461669365
$ws.Range("I11").Value = 361.699795256563
$ws.Range("B12").Value = 1999.15755542445
$ws.Range("I12").Value = 360.157555424453
$ws.Range("B13").Value = 2046.9955825325
$ws.Range("I13").Value = 228.9955825325
$ws.Range("B14").Value = 2341.43493210799
$ws.Range("I14").Value = 503.434932107989
$ws.Range("B15").Value = 2336.87402497217
$ws.Range("E15").Value = 2785.02425391867
$ws.Range("I15").Value = 487.874024972174
$ws.Range("B16").Value = 3062.72866790939
$ws.Range("C16").Value = 2434.583887596
$ws.Range("D16").Value = 2316.1680831505
$ws.Range("I16").Value = 71.7286679093918
$ws.Range("B17").Value = 2904.28289574168
$ws.Range("E17").Value = 3484.54786297648
$ws.Range("I17").Value = 107.282895741678
$ws.Range("B18").Value = 2693.09625824875
$ws.Range("E18").Value = 3469.11058030745
$ws.Range("F18").Value = 3567.55666215005
$ws.Range("I18").Value = 184.09625824875
$ws.Range("B19").Value = 2305.04720534821
$ws.Range("E19").Value = 3017.94659653821
$ws.Range("F19").Value = 3162.40986262943
$ws.Range("I19").Value = 39.0472053482094
$ws.Range("B20").Value = 2295.17546947856
$ws.Range("I20").Value = 155.17546947856
$ws.Range("B21").Value = 2258.72918208457
$ws.Range("I21").Value = 149.729182084568
$ws.Range("B22").Value = 2036.64981118202
$ws.Range("E22").Value = 2544.32579195792
$ws.Range("I22").Value = 3.64981118201649
$ws.Range("B23").Value = 1884.75774798013
$ws.Range("E23").Value = 2389.13528390233
$ws.Range("I23").Value = 38.7577479801294
$ws.Range("B24").Value = 2002.25948934241
$ws.Range("I24").Value = -52.7405106575854
$ws.Range("B25").Value = 2020.77160986269
$ws.Range("I25").Value = -348.22839013731
$ws.Range("B26").Value = 2291.62856058815
$ws.Range("I26").Value = -238.37143941185
$ws.Range("B27").Value = 2252.12702590127
$ws.Range("I27").Value = -190.872974098728
$ws.Range("B28").Value = 3135.94356684826
$ws.Range("I28").Value = 4.94356684825516
$ws.Range("B29").Value = 2785.89074918777
$ws.Range("I29").Value = 260.890749187771
$ws.Range("B30").Value = 2723.3634492245
$ws.Range("E30").Value = 3614.5359941938
$ws.Range("F30").Value = 3743.02424519787
$ws.Range("I30").Value = 220.363449224496
$ws.Range("B31").Value = 2408.58889218561
$ws.Range("E31").Value = 3354.52292458073
$ws.Range("F31").Value = 3513.55606574768
$ws.Range("I31").Value = -2.41110781439056
$ws.Range("B32").Value = 2353.78871191628
$ws.Range("I32").Value = 128.788711916285
$ws.Range("B33").Value = 2231.76644336784
$ws.Range("I33").Value = -23.2335566321572
$ws.Range("B34").Value = 2052.10171754545
$ws.Range("I34").Value = 106.101717545454
$ws.Range("B35").Value = 1881.32265579249
$ws.Range("E35").Value = 2497.10883653134
$ws.Range("I35").Value = 198.322655792489
$ws.Range("B36").Value = 2022.801690009
$ws.Range("I36").Value = 290.801690009004
$ws.Range("B37").Value = 2011.72801601862
$ws.Range("I37").Value = 824.728016018621
$ws.Range("B38").Value = 2290.83141823981
$ws.Range("I38").Value = 1146.83141823981
$ws.Range("B39").Value = 2261.25303197829
$ws.Range("I39").Value = 54.2530319782854
$ws.Range("B40").Value = 3081.48039804365
$ws.Range("I40").Value = -35.5196019563487
$ws.Range("B41").Value = 2755.12858811495
$ws.Range("I41").Value = -250.871411885052
$ws.Range("B42").Value = 2732.88082350603
$ws.Range("E42").Value = 3702.32744057989
$ws.Range("I42").Value = -89.1191764939731
$ws.Range("B43").Value = 2465.55531490631
$ws.Range("E43").Value = 3551.72043461691
$ws.Range("F43").Value = 3663.30260049759
$ws.Range("I43").Value = -63.4446850936879
$ws.Range("B44").Value = 2365.27328394649
$ws.Range("E44").Value = 3145.86398992096
$ws.Range("I44").Value = -254.726716053514
$ws.Range("B45").Value = 2263.27208721206
$ws.Range("I45").Value = -354.727912787941
$ws.Range("B46").Value = 2051.0612074951
$ws.Range("I46").Value = -498.9387925049
$ws.Range("B47").Value = 1903.63196583301
$ws.Range("I47").Value = -639.368034166995
$ws.Range("B48").Value = 2036.33958993833
$ws.Range("I48").Value = -714.660410061669
$ws.Range("B49").Value = 2020.95736926911
$ws.Range("I49").Value = -647.042630730891

$ws.Range("J22").Value = "Decrease"
$ws.Range("J28").Value = "Decrease"
